# Applies the "automatic update" diff: the per-observation data in rows
# 2,3,4,5,7,8,9 (columns A,B,D,E,F,G,H plus the optional I/M activity-count
# fields, and the Ost/Nord coordinates in Q/R) is re-shuffled among those
# rows (row 6 is left untouched). The mapping below gives, for each
# destination row, which row's original data should end up there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destRow -> sourceRow (source row refers to the *original* data, captured
# below before any writes happen)
$mapping = @{
    2 = 7
    3 = 5
    4 = 9
    5 = 4
    7 = 8
    8 = 3
    9 = 2
}

# Capture the original values for every source row before overwriting
# anything (rows are being permuted, so we must snapshot first).
$orig = @{}
foreach ($r in 2,3,4,5,7,8,9) {
    $row = @{}
    $row.A = $ws.Cells.Item($r, 1).Value2
    $row.B = $ws.Cells.Item($r, 2).Value2
    $row.D = $ws.Cells.Item($r, 4).Text
    $row.E = $ws.Cells.Item($r, 5).Value2
    $row.F = $ws.Cells.Item($r, 6).Text
    $row.G = $ws.Cells.Item($r, 7).Text
    $row.H = $ws.Cells.Item($r, 8).Text
    $row.I = $ws.Cells.Item($r, 9).Text
    $row.M = $ws.Cells.Item($r, 13).Text
    $row.Q = $ws.Cells.Item($r, 17).Value2
    $row.R = $ws.Cells.Item($r, 18).Value2
    $orig[$r] = $row
}

foreach ($destRow in 2,3,4,5,7,8,9) {
    $srcRow = $mapping[$destRow]
    $data = $orig[$srcRow]

    $ws.Cells.Item($destRow, 1).Value = $data.A
    $ws.Cells.Item($destRow, 2).Value = $data.B
    $ws.Cells.Item($destRow, 4).Value = $data.D
    $ws.Cells.Item($destRow, 5).Value = $data.E
    $ws.Cells.Item($destRow, 6).Value = $data.F
    $ws.Cells.Item($destRow, 7).Value = $data.G
    $ws.Cells.Item($destRow, 8).Value = $data.H

    $iCell = $ws.Cells.Item($destRow, 9)
    if ([string]::IsNullOrEmpty($data.I)) {
        $iCell.Value = ""
    } else {
        $iCell.NumberFormat = "@"
        $iCell.Value = $data.I
    }

    $mCell = $ws.Cells.Item($destRow, 13)
    if ([string]::IsNullOrEmpty($data.M)) {
        $mCell.Value = ""
    } else {
        $mCell.Value = $data.M
    }

    $ws.Cells.Item($destRow, 17).Value = $data.Q
    $ws.Cells.Item($destRow, 18).Value = $data.R
}
